$wb = $excel.ActiveWorkbook

# ALC row 17
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 2385182.2  # H17: 2692902.5 -> 2385182.2
$ws.Cells.Item(17, 10).Value = 2385182.2  # J17: 2692902.5 -> 2385182.2
$ws.Cells.Item(17, 12).Value = 7155546.600000001  # L17: 8078707.5 -> 7155546.600000001
$ws.Cells.Item(17, 14).Value = -7155882.600000001  # N17: -8079043.5 -> -7155882.600000001

# ALC row 21
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(21, 8).Value = 23374.875  # H21: 16079.667 -> 23374.875
$ws.Cells.Item(21, 9).Value = 2339  # I21: 1567.7142 -> 2339
$ws.Cells.Item(21, 10).Value = 35996.4  # J21: 36396.4 -> 35996.4
$ws.Cells.Item(21, 11).Value = 2339  # K21: 1567.7142 -> 2339
$ws.Cells.Item(21, 12).Value = 35996.4  # L21: 36396.4 -> 35996.4
$ws.Cells.Item(21, 13).Value = -1871  # M21: -1099.7142 -> -1871
$ws.Cells.Item(21, 14).Value = -36932.4  # N21: -37332.4 -> -36932.4

# ALC row 23
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(23, 8).Value = 23374.875  # H23: 16079.667 -> 23374.875
$ws.Cells.Item(23, 9).Value = 2339  # I23: 1567.7142 -> 2339
$ws.Cells.Item(23, 10).Value = 35996.4  # J23: 36396.4 -> 35996.4
$ws.Cells.Item(23, 11).Value = 2339  # K23: 1567.7142 -> 2339
$ws.Cells.Item(23, 12).Value = 35996.4  # L23: 36396.4 -> 35996.4
$ws.Cells.Item(23, 13).Value = -2105  # M23: -1333.7142 -> -2105
$ws.Cells.Item(23, 14).Value = -36464.4  # N23: -36864.4 -> -36464.4

# ALC row 29
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(29, 8).Value = 50  # H29: 500 -> 50
$ws.Cells.Item(29, 9).Value = 50  # I29: 500 -> 50
$ws.Cells.Item(29, 11).Value = 150  # K29: 1500 -> 150
$ws.Cells.Item(29, 13).Value = 131  # M29: -1219 -> 131

# ALC row 40
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40, 8).Value = 1688.421  # H40: 1705.5555 -> 1688.421
$ws.Cells.Item(40, 9).Value = 1854.2858  # I40: 1857.1428 -> 1854.2858
$ws.Cells.Item(40, 10).Value = 1591.6666  # J40: 1609.091 -> 1591.6666
$ws.Cells.Item(40, 11).Value = 1854.2858  # K40: 1857.1428 -> 1854.2858
$ws.Cells.Item(40, 12).Value = 1591.6666  # L40: 1609.091 -> 1591.6666
$ws.Cells.Item(40, 13).Value = -1679.2858  # M40: -1682.1428 -> -1679.2858
$ws.Cells.Item(40, 14).Value = -1941.6666  # N40: -1959.091 -> -1941.6666

# ALC row 64
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(64, 8).Value = 8550402  # H64: 3666185.8 -> 8550402
$ws.Cells.Item(64, 9).Value = 25643540  # I64: 5497428.5 -> 25643540
$ws.Cells.Item(64, 10).Value = 3833.3333  # J64: 3700 -> 3833.3333
$ws.Cells.Item(64, 11).Value = 25643540  # K64: 5497428.5 -> 25643540
$ws.Cells.Item(64, 12).Value = 3833.3333  # L64: 3700 -> 3833.3333
$ws.Cells.Item(64, 13).Value = -25643292  # M64: -5497180.5 -> -25643292
$ws.Cells.Item(64, 14).Value = -4329.3333  # N64: -4196 -> -4329.3333

# ALC row 67
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(67, 8).Value = 8550402  # H67: 3666185.8 -> 8550402
$ws.Cells.Item(67, 9).Value = 25643540  # I67: 5497428.5 -> 25643540
$ws.Cells.Item(67, 10).Value = 3833.3333  # J67: 3700 -> 3833.3333
$ws.Cells.Item(67, 11).Value = 25643540  # K67: 5497428.5 -> 25643540
$ws.Cells.Item(67, 12).Value = 3833.3333  # L67: 3700 -> 3833.3333
$ws.Cells.Item(67, 13).Value = -25642682  # M67: -5496570.5 -> -25642682
$ws.Cells.Item(67, 14).Value = -5549.3333  # N67: -5416 -> -5549.3333

# ALC row 118
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(118, 8).Value = 3180  # H118: 1741 -> 3180
$ws.Cells.Item(118, 9).Value = 0  # I118: 246.28572 -> 0
$ws.Cells.Item(118, 10).Value = 3180  # J118: 3484.8333 -> 3180
$ws.Cells.Item(118, 11).Value = 0  # K118: 738.85716 -> 0
$ws.Cells.Item(118, 12).Value = 9540  # L118: 10454.4999 -> 9540
$ws.Cells.Item(118, 13).ClearContents()  # M118: remove (was 918.14284)
$ws.Cells.Item(118, 14).Value = -12854  # N118: -13768.4999 -> -12854

# ALC row 135
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(135, 8).Value = 9435274  # H135: 13891323 -> 9435274
$ws.Cells.Item(135, 9).Value = 1053.9762  # I135: 1844 -> 1053.9762
$ws.Cells.Item(135, 10).Value = 45456840  # J135: 50003970 -> 45456840
$ws.Cells.Item(135, 11).Value = 9485.785800000001  # K135: 16596 -> 9485.785800000001
$ws.Cells.Item(135, 12).Value = 409111560  # L135: 450035730 -> 409111560
$ws.Cells.Item(135, 13).Value = -6950.785800000001  # M135: -14061 -> -6950.785800000001
$ws.Cells.Item(135, 14).Value = -409116630  # N135: -450040800 -> -409116630

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(137, 8).Value = 8696472  # H137: 15385514 -> 8696472
$ws.Cells.Item(137, 9).Value = 840.3333  # I137: 953.0909 -> 840.3333
$ws.Cells.Item(137, 11).Value = 2520.9999  # K137: 2859.2727 -> 2520.9999
$ws.Cells.Item(137, 13).Value = 29.0001000000002  # M137: -309.2727 -> 29.0001000000002

# ALC row 140
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(140, 8).Value = 79000  # H140: 98333.336 -> 79000
$ws.Cells.Item(140, 10).Value = 79000  # J140: 98333.336 -> 79000
$ws.Cells.Item(140, 12).Value = 79000  # L140: 98333.336 -> 79000
$ws.Cells.Item(140, 14).Value = -89360  # N140: -108693.336 -> -89360

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 9079.138999999999  # H32: 10579.5 -> 9079.138999999999
$ws.Cells.Item(32, 9).Value = 8440.290000000001  # I32: 9964.164000000001 -> 8440.290000000001
$ws.Cells.Item(32, 10).Value = 13040  # J32: 15414.286 -> 13040
$ws.Cells.Item(32, 11).Value = 8440.290000000001  # K32: 9964.164000000001 -> 8440.290000000001
$ws.Cells.Item(32, 12).Value = 13040  # L32: 15414.286 -> 13040
$ws.Cells.Item(32, 13).Value = -8153.290000000001  # M32: -9677.164000000001 -> -8153.290000000001
$ws.Cells.Item(32, 14).Value = -13614  # N32: -15988.286 -> -13614

# BSM row 107
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(107, 8).Value = 1152.0454  # H107: 1497.9231 -> 1152.0454
$ws.Cells.Item(107, 9).Value = 1178.6842  # I107: 1456.0834 -> 1178.6842
$ws.Cells.Item(107, 10).Value = 983.3333  # J107: 2000 -> 983.3333
$ws.Cells.Item(107, 11).Value = 1178.6842  # K107: 1456.0834 -> 1178.6842
$ws.Cells.Item(107, 12).Value = 983.3333  # L107: 2000 -> 983.3333
$ws.Cells.Item(107, 13).Value = 741.3158000000001  # M107: 463.9166 -> 741.3158000000001
$ws.Cells.Item(107, 14).Value = -4823.3333  # N107: -5840 -> -4823.3333

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 2226.9348  # H134: 2158.4329 -> 2226.9348
$ws.Cells.Item(134, 9).Value = 1732.8182  # I134: 1255.3148 -> 1732.8182
$ws.Cells.Item(134, 10).Value = 3481.2307  # J134: 5909.846 -> 3481.2307
$ws.Cells.Item(134, 11).Value = 5198.4546  # K134: 3765.9444 -> 5198.4546
$ws.Cells.Item(134, 12).Value = 10443.6921  # L134: 17729.538 -> 10443.6921
$ws.Cells.Item(134, 13).Value = -2663.4546  # M134: -1230.9444 -> -2663.4546
$ws.Cells.Item(134, 14).Value = -15513.6921  # N134: -22799.538 -> -15513.6921

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 5653247  # H31: 6670742.5 -> 5653247
$ws.Cells.Item(31, 9).Value = 3497.2449  # I31: 4228.1284 -> 3497.2449
$ws.Cells.Item(31, 10).Value = 33337020  # J31: 30306566 -> 33337020
$ws.Cells.Item(31, 11).Value = 3497.2449  # K31: 4228.1284 -> 3497.2449
$ws.Cells.Item(31, 12).Value = 33337020  # L31: 30306566 -> 33337020
$ws.Cells.Item(31, 13).Value = -3202.2449  # M31: -3933.1284 -> -3202.2449
$ws.Cells.Item(31, 14).Value = -33337610  # N31: -30307156 -> -33337610

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(34, 8).Value = 5653247  # H34: 6670742.5 -> 5653247
$ws.Cells.Item(34, 9).Value = 3497.2449  # I34: 4228.1284 -> 3497.2449
$ws.Cells.Item(34, 10).Value = 33337020  # J34: 30306566 -> 33337020
$ws.Cells.Item(34, 11).Value = 3497.2449  # K34: 4228.1284 -> 3497.2449
$ws.Cells.Item(34, 12).Value = 33337020  # L34: 30306566 -> 33337020
$ws.Cells.Item(34, 13).Value = -3295.2449  # M34: -4026.1284 -> -3295.2449
$ws.Cells.Item(34, 14).Value = -33337424  # N34: -30306970 -> -33337424

# CRP row 62
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(62, 8).Value = 0  # H62: 2150 -> 0
$ws.Cells.Item(62, 9).Value = 0  # I62: 2200 -> 0
$ws.Cells.Item(62, 10).Value = 0  # J62: 2060 -> 0
$ws.Cells.Item(62, 11).Value = 0  # K62: 2200 -> 0
$ws.Cells.Item(62, 12).Value = 0  # L62: 2060 -> 0
$ws.Cells.Item(62, 13).ClearContents()  # M62: remove (was -1576)
$ws.Cells.Item(62, 14).ClearContents()  # N62: remove (was -3308)

# CRP row 65
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(65, 8).Value = 0  # H65: 2150 -> 0
$ws.Cells.Item(65, 9).Value = 0  # I65: 2200 -> 0
$ws.Cells.Item(65, 10).Value = 0  # J65: 2060 -> 0
$ws.Cells.Item(65, 11).Value = 0  # K65: 11000 -> 0
$ws.Cells.Item(65, 12).Value = 0  # L65: 10300 -> 0
$ws.Cells.Item(65, 13).ClearContents()  # M65: remove (was -7880)
$ws.Cells.Item(65, 14).ClearContents()  # N65: remove (was -16540)

# CRP row 105
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(105, 8).Value = 1933.75  # H105: 2046.6666 -> 1933.75
$ws.Cells.Item(105, 9).Value = 2203.3333  # I105: 2152.5 -> 2203.3333
$ws.Cells.Item(105, 10).Value = 1125  # J105: 1200 -> 1125
$ws.Cells.Item(105, 11).Value = 2203.3333  # K105: 2152.5 -> 2203.3333
$ws.Cells.Item(105, 12).Value = 1125  # L105: 1200 -> 1125
$ws.Cells.Item(105, 13).Value = -456.3332999999998  # M105: -405.5 -> -456.3332999999998
$ws.Cells.Item(105, 14).Value = -4619  # N105: -4694 -> -4619

# CUL row 125
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(125, 8).Value = 3229.2307  # H125: 3016.3635 -> 3229.2307
$ws.Cells.Item(125, 9).Value = 2000  # I125: 2066.6667 -> 2000
$ws.Cells.Item(125, 10).Value = 3598  # J125: 3372.5 -> 3598
$ws.Cells.Item(125, 11).Value = 6000  # K125: 6200.000100000001 -> 6000
$ws.Cells.Item(125, 12).Value = 10794  # L125: 10117.5 -> 10794
$ws.Cells.Item(125, 13).Value = -1080  # M125: -1280.000100000001 -> -1080
$ws.Cells.Item(125, 14).Value = -20634  # N125: -19957.5 -> -20634

# CUL row 137
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(137, 8).Value = 9266063  # H137: 6708.35 -> 9266063
$ws.Cells.Item(137, 9).Value = 23813838  # I137: 3041.75 -> 23813838
$ws.Cells.Item(137, 10).Value = 8387.817999999999  # J137: 12208.25 -> 8387.817999999999
$ws.Cells.Item(137, 11).Value = 71441514  # K137: 9125.25 -> 71441514
$ws.Cells.Item(137, 12).Value = 25163.454  # L137: 36624.75 -> 25163.454
$ws.Cells.Item(137, 13).Value = -71436414  # M137: -4025.25 -> -71436414
$ws.Cells.Item(137, 14).Value = -35363.454  # N137: -46824.75 -> -35363.454

# GSM row 102
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 4407.6  # H102: 5272 -> 4407.6
$ws.Cells.Item(102, 9).Value = 4572.8237  # I102: 5687.5386 -> 4572.8237
$ws.Cells.Item(102, 11).Value = 4572.8237  # K102: 5687.5386 -> 4572.8237
$ws.Cells.Item(102, 13).Value = -2950.8237  # M102: -4065.5386 -> -2950.8237

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 3076.2144  # H132: 3041.7544 -> 3076.2144
$ws.Cells.Item(132, 9).Value = 2416.0852  # I132: 2559.3953 -> 2416.0852
$ws.Cells.Item(132, 10).Value = 6523.5557  # J132: 4523.2856 -> 6523.5557
$ws.Cells.Item(132, 11).Value = 7248.2556  # K132: 7678.1859 -> 7248.2556
$ws.Cells.Item(132, 12).Value = 19570.6671  # L132: 13569.8568 -> 19570.6671
$ws.Cells.Item(132, 13).Value = -4718.2556  # M132: -5148.1859 -> -4718.2556
$ws.Cells.Item(132, 14).Value = -24630.6671  # N132: -18629.8568 -> -24630.6671

# LTW row 7
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 5919.2583  # H7: 5148.2 -> 5919.2583
$ws.Cells.Item(7, 9).Value = 7669.2  # I7: 5905.467 -> 7669.2
$ws.Cells.Item(7, 10).Value = 5085.952  # J7: 4580.25 -> 5085.952
$ws.Cells.Item(7, 11).Value = 7669.2  # K7: 5905.467 -> 7669.2
$ws.Cells.Item(7, 12).Value = 5085.952  # L7: 4580.25 -> 5085.952
$ws.Cells.Item(7, 13).Value = -7557.2  # M7: -5793.467 -> -7557.2
$ws.Cells.Item(7, 14).Value = -5309.952  # N7: -4804.25 -> -5309.952

# LTW row 40
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 11257.714  # H40: 14600.8 -> 11257.714
$ws.Cells.Item(40, 9).Value = 11960.8  # I40: 25502 -> 11960.8
$ws.Cells.Item(40, 10).Value = 9500  # J40: 7333.3335 -> 9500
$ws.Cells.Item(40, 11).Value = 11960.8  # K40: 25502 -> 11960.8
$ws.Cells.Item(40, 12).Value = 9500  # L40: 7333.3335 -> 9500
$ws.Cells.Item(40, 13).Value = -11824.8  # M40: -25366 -> -11824.8
$ws.Cells.Item(40, 14).Value = -9772  # N40: -7605.3335 -> -9772

# LTW row 61
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(61, 8).Value = 1725.4  # H61: 1496.4 -> 1725.4
$ws.Cells.Item(61, 9).Value = 1563  # I61: 1366.2 -> 1563
$ws.Cells.Item(61, 10).Value = 2375  # J61: 1887 -> 2375
$ws.Cells.Item(61, 11).Value = 1563  # K61: 1366.2 -> 1563
$ws.Cells.Item(61, 12).Value = 2375  # L61: 1887 -> 2375
$ws.Cells.Item(61, 13).Value = -1361  # M61: -1164.2 -> -1361
$ws.Cells.Item(61, 14).Value = -2779  # N61: -2291 -> -2779

# LTW row 94
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(94, 8).Value = 41331.668  # H94: 41331.5 -> 41331.668
$ws.Cells.Item(94, 10).Value = 41331.668  # J94: 41331.5 -> 41331.668
$ws.Cells.Item(94, 12).Value = 41331.668  # L94: 41331.5 -> 41331.668
$ws.Cells.Item(94, 14).Value = -42683.668  # N94: -42683.5 -> -42683.668

# LTW row 113
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(113, 8).Value = 1725.4  # H113: 1496.4 -> 1725.4
$ws.Cells.Item(113, 9).Value = 1563  # I113: 1366.2 -> 1563
$ws.Cells.Item(113, 10).Value = 2375  # J113: 1887 -> 2375
$ws.Cells.Item(113, 11).Value = 1563  # K113: 1366.2 -> 1563
$ws.Cells.Item(113, 12).Value = 2375  # L113: 1887 -> 2375
$ws.Cells.Item(113, 13).Value = 607  # M113: 803.8 -> 607
$ws.Cells.Item(113, 14).Value = -6715  # N113: -6227 -> -6715

# LTW row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(122, 8).Value = 4846.9653  # H122: 4431.6113 -> 4846.9653
$ws.Cells.Item(122, 9).Value = 4867.913  # I122: 4590.346 -> 4867.913
$ws.Cells.Item(122, 10).Value = 4766.6665  # J122: 4018.9 -> 4766.6665
$ws.Cells.Item(122, 11).Value = 14603.739  # K122: 13771.038 -> 14603.739
$ws.Cells.Item(122, 12).Value = 14299.9995  # L122: 12056.7 -> 14299.9995
$ws.Cells.Item(122, 13).Value = -12153.739  # M122: -11321.038 -> -12153.739
$ws.Cells.Item(122, 14).Value = -19199.9995  # N122: -16956.7 -> -19199.9995

# LTW row 126
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(126, 8).Value = 5919.2583  # H126: 5148.2 -> 5919.2583
$ws.Cells.Item(126, 9).Value = 7669.2  # I126: 5905.467 -> 7669.2
$ws.Cells.Item(126, 10).Value = 5085.952  # J126: 4580.25 -> 5085.952
$ws.Cells.Item(126, 11).Value = 23007.6  # K126: 17716.401 -> 23007.6
$ws.Cells.Item(126, 12).Value = 15257.856  # L126: 13740.75 -> 15257.856
$ws.Cells.Item(126, 13).Value = -20537.6  # M126: -15246.401 -> -20537.6
$ws.Cells.Item(126, 14).Value = -20197.856  # N126: -18680.75 -> -20197.856

# WVR row 28
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(28, 8).Value = 5016.5  # H28: 4764.25 -> 5016.5
$ws.Cells.Item(28, 9).Value = 0  # I28: 4000 -> 0
$ws.Cells.Item(28, 10).Value = 5016.5  # J28: 5019 -> 5016.5
$ws.Cells.Item(28, 11).Value = 0  # K28: 4000 -> 0
$ws.Cells.Item(28, 12).Value = 5016.5  # L28: 5019 -> 5016.5
$ws.Cells.Item(28, 13).ClearContents()  # M28: remove (was -3652)
$ws.Cells.Item(28, 14).Value = -5712.5  # N28: -5715 -> -5712.5

# WVR row 62
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 11975.25  # H62: 10280.667 -> 11975.25
$ws.Cells.Item(62, 9).Value = 6300  # I62: 6134.9 -> 6300
$ws.Cells.Item(62, 10).Value = 16029  # J62: 15462.875 -> 16029
$ws.Cells.Item(62, 11).Value = 6300  # K62: 6134.9 -> 6300
$ws.Cells.Item(62, 12).Value = 16029  # L62: 15462.875 -> 16029
$ws.Cells.Item(62, 13).Value = -5676  # M62: -5510.9 -> -5676
$ws.Cells.Item(62, 14).Value = -17277  # N62: -16710.875 -> -17277

# WVR row 65
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(65, 8).Value = 11975.25  # H65: 10280.667 -> 11975.25
$ws.Cells.Item(65, 9).Value = 6300  # I65: 6134.9 -> 6300
$ws.Cells.Item(65, 10).Value = 16029  # J65: 15462.875 -> 16029
$ws.Cells.Item(65, 11).Value = 31500  # K65: 30674.5 -> 31500
$ws.Cells.Item(65, 12).Value = 80145  # L65: 77314.375 -> 80145
$ws.Cells.Item(65, 13).Value = -28380  # M65: -27554.5 -> -28380
$ws.Cells.Item(65, 14).Value = -86385  # N65: -83554.375 -> -86385

# WVR row 81
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 1420  # H81: 1053.3334 -> 1420
$ws.Cells.Item(81, 9).Value = 950  # I81: 433.33334 -> 950
$ws.Cells.Item(81, 10).Value = 1733.3334  # J81: 1260 -> 1733.3334
$ws.Cells.Item(81, 11).Value = 1900  # K81: 866.66668 -> 1900
$ws.Cells.Item(81, 12).Value = 3466.6668  # L81: 2520 -> 3466.6668
$ws.Cells.Item(81, 13).Value = -839  # M81: 194.33332 -> -839
$ws.Cells.Item(81, 14).Value = -5588.6668  # N81: -4642 -> -5588.6668

# WVR row 84
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(84, 8).Value = 1420  # H84: 1053.3334 -> 1420
$ws.Cells.Item(84, 9).Value = 950  # I84: 433.33334 -> 950
$ws.Cells.Item(84, 10).Value = 1733.3334  # J84: 1260 -> 1733.3334
$ws.Cells.Item(84, 11).Value = 9500  # K84: 4333.3334 -> 9500
$ws.Cells.Item(84, 12).Value = 17333.334  # L84: 12600 -> 17333.334
$ws.Cells.Item(84, 13).Value = -4196  # M84: 970.6665999999996 -> -4196
$ws.Cells.Item(84, 14).Value = -27941.334  # N84: -23208 -> -27941.334
